# case 1 data update: row 1 grows from 14 values (10 counts + 4 fractions)
# to 17 values (12 counts + 5 fractions), and several columns get wider.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row 1 values --------------------------------------------------
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 31
$ws.Range("D1").Value = 13
$ws.Range("E1").Value = 30
$ws.Range("F1").Value = 32
$ws.Range("G1").Value = 22
$ws.Range("H1").Value = 15
$ws.Range("I1").Value = 33
$ws.Range("J1").Value = 23
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 21
$ws.Range("M1").Value = 0.08399999999999999
$ws.Range("N1").Value = 0.031
$ws.Range("O1").Value = 0.072
$ws.Range("P1").Value = 0.012
$ws.Range("Q1").Value = 0.066

# --- column widths -------------------------------------------------------
# Columns A-D and M already have the right width and are left untouched.
# Columns E-L widen to match column C/D's width; N-Q widen to match M's.
# NOTE: Excel's ColumnWidth setter snaps to a pixel grid, so the requested
# "characters" width is offset by the fixed 5/6 pixel->char conversion
# before assignment; this lands on the closest attainable grid value
# (3.1666... and 5.6666... respectively) to the canonical 3.140625 /
# 5.7109375 widths stored in the workbook.
$ws.Columns.Item(5).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(6).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(7).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(8).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(9).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(10).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(11).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(12).ColumnWidth = 2.3072916666666665
$ws.Columns.Item(14).ColumnWidth = 4.877604166666667
$ws.Columns.Item(15).ColumnWidth = 4.877604166666667
$ws.Columns.Item(16).ColumnWidth = 4.877604166666667
$ws.Columns.Item(17).ColumnWidth = 4.877604166666667
